$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Name and Department columns for rows 2-5
$ws.Range("B2").Value = "Rahul"
$ws.Range("C2").Value = "QA"

$ws.Range("B3").Value = "Nitin"
$ws.Range("C3").Value = "QA"

$ws.Range("B4").Value = "Binu"
$ws.Range("C4").Value = "Dev"

$ws.Range("B5").Value = "joy"
$ws.Range("C5").Value = "Dev"

# Remove row 6 entirely (was Employee 5 / Mathew / Bussiness)
$ws.Rows("6:6").Delete()
